# Mark "OK" in column E for the checklist rows that don't need the
# command factory pattern (i.e. don't need undo/redo command support):
#   row 5  -> "3. Show current player"
#   row 6  -> "4. display all player"
#   row 7  -> "5. set current player by playerID"
#   row 11 -> "9. show undo/ redo list"
#   row 12 -> "10. undo"
#   row 13 -> "11. redo"
#   row 14 -> "12. exit"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5, 6, 7, 11, 12, 13, 14)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "OK"
}

$ws.Range("E8").Select()
